# Fix formatting bug introduced when scrapping floating point numbers:
# a naive scraper had swapped the decimal separator for a comma (and left
# thousands separators as dots). Restore plain numeric-style text by
# stripping the "." thousands separators and turning the decimal "," into
# a ".". The same substitution incidentally also needs to be undone on a
# few free-text "Razon social"/"Nombre Fantasia" cells that happened to
# contain a literal comma (e.g. "LASTNAME A, LASTNAME B") - so the fix is
# applied uniformly to every text cell that still contains a comma.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -is [string] -and $val.Contains(",")) {
            $newVal = $val.Replace(".", "").Replace(",", ".")

            # Writing a plain-looking decimal string (e.g. "3870.00") back
            # through .Value would be auto-coerced to the NUMBER 3870,
            # losing the original "text" cell type the source file has
            # (t="s" in the sst, no number format). Force text entry the
            # same way Excel does for an explicitly text-formatted cell,
            # then restore the cell to the plain "Normal" style so its
            # format stays exactly as it was (General, no explicit style).
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        }
    }
}
